# Updated cryptos list with latest prices and volume changes
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'23.338.37"
$ws.Range("E2").Value = "  -0.39%  "

# Row 3
$ws.Range("D3").Value = "'1.625.67"
$ws.Range("E3").Value = "  -0.78%  "

# Row 4
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  +0.01%  "

# Row 5
$ws.Range("D5").Value = "'1.001"
$ws.Range("E5").Value = "  +0.04%  "

# Row 6
$ws.Range("D6").Value = "'303.29"
$ws.Range("E6").Value = "  -0.53%  "

# Row 7
$ws.Range("D7").Value = "'0.3739"
$ws.Range("E7").Value = "  +0.04%  "

# Row 8
$ws.Range("B8").Value = "Cardano"
$ws.Range("C8").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D8").Value = "'0.3623"
$ws.Range("E8").Value = "  +0.12%  "

# Row 9
$ws.Range("B9").Value = "OKB"
$ws.Range("C9").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D9").Value = "'51.40"
$ws.Range("E9").Value = "  -0.88%  "

# Row 10
$ws.Range("D10").Value = "'0.08143"
$ws.Range("E10").Value = "  +0.32%  "

# Row 11
$ws.Range("D11").Value = "'1.225"
$ws.Range("E11").Value = "  -2.38%  "

# Row 12
$ws.Range("D12").Value = "'1.001"

# Row 13
$ws.Range("D13").Value = "'22.26"
$ws.Range("E13").Value = "  -2.50%  "

# Row 14
$ws.Range("D14").Value = "'6.462"
$ws.Range("E14").Value = "  -1.97%  "

# Row 15
$ws.Range("D15").Value = "'0.00001240"
$ws.Range("E15").Value = "  -2.11%  "

# Row 16
$ws.Range("D16").Value = "'7.279"
$ws.Range("E16").Value = "  -0.08%  "

# Row 17
$ws.Range("D17").Value = "'1.625.94"
$ws.Range("E17").Value = "  -0.73%  "

# Row 18
$ws.Range("D18").Value = "'93.85"
$ws.Range("E18").Value = "  -0.37%  "

# Row 19
$ws.Range("D19").Value = "'0.06945"
$ws.Range("E19").Value = "  +0.62%  "

# Row 20
$ws.Range("D20").Value = "'17.50"
$ws.Range("E20").Value = "  -3.37%  "

# Row 21
$ws.Range("D21").Value = "'6.531"
$ws.Range("E21").Value = "  +0.52%  "

# Row 22
$ws.Range("D22").Value = "'1.001"
$ws.Range("E22").Value = "  +0.09%  "

# Row 23
$ws.Range("D23").Value = "'12.52"
$ws.Range("E23").Value = "  -1.56%  "

# Row 24
$ws.Range("D24").Value = "'23.342.75"
$ws.Range("E24").Value = "  -0.41%  "

# Row 25
$ws.Range("D25").Value = "'2.464"
$ws.Range("E25").Value = "  +1.78%  "

# Row 26
$ws.Range("D26").Value = "'3.098"
$ws.Range("E26").Value = "  +2.12%  "

# Row 27
$ws.Range("D27").Value = "'21.15"
$ws.Range("E27").Value = "  -0.19%  "

# Row 28
$ws.Range("D28").Value = "'150.17"
$ws.Range("E28").Value = "  -0.92%  "

# Row 29
$ws.Range("D29").Value = "'5.259"
$ws.Range("E29").Value = "  -1.31%  "

# Row 30
$ws.Range("D30").Value = "'132.54"
$ws.Range("E30").Value = "  -2.18%  "

# Row 31
$ws.Range("D31").Value = "'1.794.91"
$ws.Range("E31").Value = "  -1.33%  "

# Row 32
$ws.Range("D32").Value = "'6.712"
$ws.Range("E32").Value = "  -0.23%  "

# Row 33
$ws.Range("D33").Value = "'2.160"
$ws.Range("E33").Value = "  -4.90%  "

# Row 34
$ws.Range("D34").Value = "'1.044"
$ws.Range("E34").Value = "  +9.32%  "

# Row 35
$ws.Range("D35").Value = "'10.86"
$ws.Range("E35").Value = "  +5.99%  "

# Row 36
$ws.Range("D36").Value = "'0.02755"
$ws.Range("E36").Value = "  -2.07%  "

# Row 37
$ws.Range("B37").Value = "Stellar"
$ws.Range("C37").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D37").Value = "'0.08773"
$ws.Range("E37").Value = "  -0.05%  "

# Row 38
$ws.Range("B38").Value = "Algorand"
$ws.Range("C38").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D38").Value = "'0.2490"
$ws.Range("E38").Value = "  -1.01%  "

# Row 39
$ws.Range("D39").Value = "'0.07089"
$ws.Range("E39").Value = "  -2.13%  "

# Row 40
$ws.Range("D40").Value = "'5.977"
$ws.Range("E40").Value = "  -1.29%  "

# Row 41
$ws.Range("B41").Value = "TrustWalletToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D41").Value = "'1.338"
$ws.Range("E41").Value = "  -2.46%  "

# Row 42
$ws.Range("B42").Value = "TheSandbox"
$ws.Range("C42").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D42").Value = "'0.6969"
$ws.Range("E42").Value = "  -0.97%  "

# Row 43
$ws.Range("D43").Value = "'15.92"
$ws.Range("E43").Value = "  -0.79%  "

# Row 44
$ws.Range("D44").Value = "'12.07"
$ws.Range("E44").Value = "  -2.74%  "

# Row 45
$ws.Range("D45").Value = "'0.6460"
$ws.Range("E45").Value = "  -0.47%  "

# Row 46
$ws.Range("E46").Value = "  +0.06%  "

# Row 47
$ws.Range("D47").Value = "'3.962"
$ws.Range("E47").Value = "  -1.14%  "

# Row 48
$ws.Range("D48").Value = "'2.262"
$ws.Range("E48").Value = "  -2.62%  "

# Row 49
$ws.Range("D49").Value = "'0.07970"
$ws.Range("E49").Value = "  +0.04%  "

# Row 50
$ws.Range("D50").Value = "'125.85"
$ws.Range("E50").Value = "  -1.74%  "

# Row 51
$ws.Range("D51").Value = "'1.183"
$ws.Range("E51").Value = "  -1.40%  "
